# Update "Mexico Liga MX Femenil" worksheet:
#  - Several match rows had their HomeTeam/AwayTeam (and all associated
#    odds/result data) rows reordered (pairwise swaps, and one 3-way
#    rotation). We swap the data (columns B:AC) between the affected row
#    numbers, leaving column A (the running index) untouched.
#  - A new match row (row 241) is appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($sheet, $rowA, $rowB)
    $valsA = $sheet.Range("B$rowA`:AC$rowA").Value2
    $valsB = $sheet.Range("B$rowB`:AC$rowB").Value2
    $sheet.Range("B$rowA`:AC$rowA").Value2 = $valsB
    $sheet.Range("B$rowB`:AC$rowB").Value2 = $valsA
}

# Pairwise swaps of (HomeTeam/AwayTeam + odds) between two rows.
Swap-Rows $ws 47 48
Swap-Rows $ws 101 102
Swap-Rows $ws 131 132
Swap-Rows $ws 133 134
Swap-Rows $ws 149 150
Swap-Rows $ws 215 216

# Three-way rotation: new229 = old231, new230 = old229, new231 = old230
$v229 = $ws.Range("B229:AC229").Value2
$v230 = $ws.Range("B230:AC230").Value2
$v231 = $ws.Range("B231:AC231").Value2
$ws.Range("B229:AC229").Value2 = $v231
$ws.Range("B230:AC230").Value2 = $v229
$ws.Range("B231:AC231").Value2 = $v230

# Append a new row 241 with a new match.
# Copy number formatting (styles) from row 240's A and E cells so the new
# row matches the existing formatting (bold/border index column, date
# column format) without introducing new style entries.
$ws.Cells.Item(240, 1).Copy()
$ws.Cells.Item(241, 1).PasteSpecial(-4122)
$ws.Cells.Item(240, 5).Copy()
$ws.Cells.Item(241, 5).PasteSpecial(-4122)

$ws.Cells.Item(241, 1).Value = 239
$ws.Cells.Item(241, 2).Value = 7645717
$ws.Cells.Item(241, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(241, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(241, 5).Value = 45353.625
$ws.Cells.Item(241, 6).Value = "Unam Pumas Women"
$ws.Cells.Item(241, 7).Value = "Tijuana Women"
$ws.Cells.Item(241, 11).Value = 2.1
$ws.Cells.Item(241, 12).Value = 3.6
$ws.Cells.Item(241, 13).Value = 2.875
$ws.Cells.Item(241, 14).Value = 2.1
$ws.Cells.Item(241, 15).Value = 3.5
$ws.Cells.Item(241, 16).Value = 2.875
$ws.Cells.Item(241, 17).Value = -0.25
$ws.Cells.Item(241, 18).Value = 1.9
$ws.Cells.Item(241, 19).Value = 1.9
$ws.Cells.Item(241, 20).Value = 3
$ws.Cells.Item(241, 21).Value = 1.775
$ws.Cells.Item(241, 22).Value = 2.025
$ws.Cells.Item(241, 23).Value = 0
$ws.Cells.Item(241, 24).Value = 0
$ws.Cells.Item(241, 25).Value = 0
$ws.Cells.Item(241, 26).Value = 0
$ws.Cells.Item(241, 27).Value = 0
